$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 53; this shifts rows 53:161 down to 54:162
# and extends the used range / dimension to row 162 automatically.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new weekly record.
$ws.Range("A53").Value = 6
$ws.Range("B53").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C53").Value = "Metropolitana"
$ws.Range("D53").Value = 44579
$ws.Range("E53").Value = 13
$ws.Range("F53").Value = 100112001
$ws.Range("G53").Value = "Berenjena"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 200
$ws.Range("K53").Value = 10000
$ws.Range("L53").Value = 12000
$ws.Range("M53").Value = 11200
$ws.Range("N53").Value = "$/caja 60 unidades"
$ws.Range("O53").Value = "Región Metropolitana"
$ws.Range("P53").Value = 187
$ws.Range("Q53").Value = 60
$ws.Range("R53").Value = "Hortaliza"
